$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: the existing row 3 (week of 2021-10-07) becomes historical data that
# needs to be preserved on a new row 4 before row 3 is overwritten with the
# new weekly values.
$ws.Range("A4:R4").Value = $ws.Range("A3:R3").Value()
$ws.Range("D4").NumberFormat = $ws.Range("D3").NumberFormat()

# Step 2: update row 3 in place with this week's new values.
$ws.Range("D3").Value = 44503
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 950
$ws.Range("P3").Value = 950
